# Add new "Diffusivity_Pressure" column (G) and refresh the pressure
# calculation results produced by the updated PressureCalculationModel
# methods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header -----------------------------------------------
# Copy F1's formatting (bold, centered, bordered header style) onto G1,
# then set its text.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Diffusivity_Pressure"

# --- Updated per-well results ------------------------------------------
# Columns: B=Initial_Pressure, C=Calculated_Pressure, D=Difference,
#          E=Adjusted_Pressure, F=Boundary_Applied, G=Diffusivity_Pressure

$data = @(
    @{ Row = 2;  B = 246.764970144766;  C = 181.5150453311255; D = 65.24992481364052;  E = 231.764970144766;  F = $true;  G = 248.3497004189754 },
    @{ Row = 3;  B = 217.0336452784793; C = 199.1920685601137; D = 17.84157671836564;  E = 202.0336452784793; F = $true;  G = 248.3497004189754 },
    @{ Row = 4;  B = 219.5933175347968; C = 191.1166021712214; D = 28.47671536357535;  E = 204.5933175347968; F = $true;  G = 248.3497004189754 },
    @{ Row = 5;  B = 206.4928996373361; C = 201.1708477223052; D = 5.322051915030869;  E = 201.1708477223052; F = $false; G = 248.3497004189754 },
    @{ Row = 6;  B = 209.0464864939379; C = 180.0241806025757; D = 29.02230589136221;  E = 194.0464864939379; F = $true;  G = 248.3497004189754 },
    @{ Row = 7;  B = 237.1332612683818; C = 217.8679032848748; D = 19.26535798350699;  E = 222.1332612683818; F = $true;  G = 248.3497004189754 },
    @{ Row = 8;  B = 231.8836493314376; C = 192.829458799071;  D = 39.05419053236662;  E = 216.8836493314376; F = $true;  G = 248.3497004189754 },
    @{ Row = 9;  B = 219.5465280520848; C = 187.6473592943935; D = 31.89916875769134;  E = 204.5465280520848; F = $true;  G = 248.3497004189754 },
    @{ Row = 10; B = 208.3029323774192; C = 211.4293246328534; D = -3.126392255434126; E = 211.4293246328534; F = $false; G = 248.3497004189754 },
    @{ Row = 11; B = 239.3385851314313; C = 223.350607012967;  D = 15.98797811846435;  E = 224.3385851314313; F = $true;  G = 248.3497004189754 }
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
}
